$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("data2 weibull")
$ws.Cells.Item(2, 5).Value = -2.39500202175584
$ws.Cells.Item(2, 6).Value = 0.10812372759168
$ws.Cells.Item(2, 7).Value = 0.103772517733394
$ws.Cells.Item(2, 8).Value = 0.0769739957034721
$ws.Cells.Item(2, 9).Value = 0.0116907404683198
$ws.Cells.Item(2, 10).Value = 0.00592499601455814
$ws.Cells.Item(2, 11).Value = -0.00402251489432989
$ws.Cells.Item(3, 5).Value = -2.57824772385901
$ws.Cells.Item(3, 6).Value = 0.125256517287666
$ws.Cells.Item(3, 7).Value = 0.0468565894924138
$ws.Cells.Item(3, 8).Value = 0.05928365352867
$ws.Cells.Item(3, 9).Value = 0.0156891951230354
$ws.Cells.Item(3, 10).Value = 0.00351455157570738
$ws.Cells.Item(3, 11).Value = -0.00425682222495584
$ws.Cells.Item(4, 5).Value = -3.13810935882744
$ws.Cells.Item(4, 6).Value = 0.188580141371055
$ws.Cells.Item(4, 7).Value = 0.150618291702876
$ws.Cells.Item(4, 8).Value = 0.122038959982672
$ws.Cells.Item(4, 9).Value = 0.0355624697195271
$ws.Cells.Item(4, 10).Value = 0.0148935077536523
$ws.Cells.Item(4, 11).Value = -0.0152701894491901
$ws.Cells.Item(5, 5).Value = -3.58822680878078
$ws.Cells.Item(5, 6).Value = 0.250481693413968
$ws.Cells.Item(5, 7).Value = 0.327498708634429
$ws.Cells.Item(5, 8).Value = 0.118770393218704
$ws.Cells.Item(5, 9).Value = 0.0627410787355289
$ws.Cells.Item(5, 10).Value = 0.0141064063053257
$ws.Cells.Item(5, 11).Value = -0.0232967607253938
$ws.Cells.Item(6, 5).Value = -2.48739851639861
$ws.Cells.Item(6, 6).Value = 0.119237197740365
$ws.Cells.Item(6, 7).Value = -0.0173439251044164
$ws.Cells.Item(6, 8).Value = 0.0858621209145547
$ws.Cells.Item(6, 9).Value = 0.0142175093249749
$ws.Cells.Item(6, 10).Value = 0.00737230380794561
$ws.Cells.Item(6, 11).Value = -0.00675677729281162
$ws.Cells.Item(7, 5).Value = -2.61351770303824
$ws.Cells.Item(7, 6).Value = 0.0912646059721666
$ws.Cells.Item(7, 7).Value = -0.0447916602245716
$ws.Cells.Item(7, 8).Value = 0.0750187544228287
$ws.Cells.Item(7, 9).Value = 0.00832922830325483
$ws.Cells.Item(7, 10).Value = 0.00562781351515267
$ws.Cells.Item(7, 11).Value = -0.00329069460677465
$ws.Cells.Item(8, 5).Value = -2.35992376256339
$ws.Cells.Item(8, 6).Value = 0.309719915328088
$ws.Cells.Item(8, 7).Value = -0.153144231507661
$ws.Cells.Item(8, 8).Value = 0.150161873039823
$ws.Cells.Item(8, 9).Value = 0.0959264259508379
$ws.Cells.Item(8, 10).Value = 0.0225485881148278
$ws.Cells.Item(8, 11).Value = -0.0218092681903502
$ws.Cells.Item(9, 5).Value = -3.00978071494291
$ws.Cells.Item(9, 6).Value = 0.294620451678279
$ws.Cells.Item(9, 7).Value = 0.138978900167005
$ws.Cells.Item(9, 8).Value = 0.170187449721332
$ws.Cells.Item(9, 9).Value = 0.086801210547113
$ws.Cells.Item(9, 10).Value = 0.0289637680426511
$ws.Cells.Item(9, 11).Value = -0.0421769037283446
$ws.Cells.Item(10, 5).Value = -2.14218267845162
$ws.Cells.Item(10, 6).Value = 0.273912315481934
$ws.Cells.Item(10, 7).Value = 0.1969390943671
$ws.Cells.Item(10, 8).Value = 0.161633137077179
$ws.Cells.Item(10, 9).Value = 0.0750279565726748
$ws.Cells.Item(10, 10).Value = 0.0261252710014103
$ws.Cells.Item(10, 11).Value = -0.0311812098035073
$ws.Cells.Item(11, 5).Value = -2.82367267232602
$ws.Cells.Item(11, 6).Value = 0.324682196980101
$ws.Cells.Item(11, 7).Value = 0.175621856814685
$ws.Cells.Item(11, 8).Value = 0.13578283802204
$ws.Cells.Item(11, 9).Value = 0.105418529035825
$ws.Cells.Item(11, 10).Value = 0.0184369791013194
$ws.Cells.Item(11, 11).Value = -0.0377651881515243
$ws.Cells.Item(12, 5).Value = -2.85084742220134
$ws.Cells.Item(12, 6).Value = 0.373118632041817
$ws.Cells.Item(12, 7).Value = 0.191139707450786
$ws.Cells.Item(12, 8).Value = 0.184639301058125
$ws.Cells.Item(12, 9).Value = 0.139217513576757
$ws.Cells.Item(12, 10).Value = 0.034091671495233
$ws.Cells.Item(12, 11).Value = -0.0563951175806997
$ws.Cells.Item(13, 5).Value = -4.00779465950803
$ws.Cells.Item(13, 6).Value = 0.628093159789506
$ws.Cells.Item(13, 7).Value = 0.375222007563247
$ws.Cells.Item(13, 8).Value = 0.244402793342619
$ws.Cells.Item(13, 9).Value = 0.394501017374366
$ws.Cells.Item(13, 10).Value = 0.059732725393675
$ws.Cells.Item(13, 11).Value = -0.14457190032534
$ws.Cells.Item(14, 5).Value = -2.84402157787172
$ws.Cells.Item(14, 6).Value = 0.166400404319523
$ws.Cells.Item(14, 7).Value = -0.0114217344621236
$ws.Cells.Item(14, 8).Value = 0.0918101730715185
$ws.Cells.Item(14, 9).Value = 0.0276890945577009
$ws.Cells.Item(14, 10).Value = 0.00842910787942218
$ws.Cells.Item(14, 11).Value = -0.00936032305604423
$ws.Cells.Item(15, 5).Value = -3.01920941114309
$ws.Cells.Item(15, 6).Value = 0.1957543823296
$ws.Cells.Item(15, 7).Value = 0.0688710722139505
$ws.Cells.Item(15, 8).Value = 0.106918285891009
$ws.Cells.Item(15, 9).Value = 0.0383197782012431
$ws.Cells.Item(15, 10).Value = 0.0114315198578715
$ws.Cells.Item(15, 11).Value = -0.0153403431417579
$ws.Cells.Item(16, 5).Value = -2.8473150791631
$ws.Cells.Item(16, 6).Value = 0.144695958439284
$ws.Cells.Item(16, 7).Value = -0.0750160783595377
$ws.Cells.Item(16, 8).Value = 0.0858304841970044
$ws.Cells.Item(16, 9).Value = 0.0209369203886629
$ws.Cells.Item(16, 10).Value = 0.00736687201749222
$ws.Cells.Item(16, 11).Value = -0.00586329307363494
$ws.Cells.Item(17, 5).Value = -1.83658401006587
$ws.Cells.Item(17, 6).Value = 0.0953889499198419
$ws.Cells.Item(17, 7).Value = -0.197722113791001
$ws.Cells.Item(17, 8).Value = 0.0419833442232253
$ws.Cells.Item(17, 9).Value = 0.00909905176681011
$ws.Cells.Item(17, 10).Value = 0.00176260119216583
$ws.Cells.Item(17, 11).Value = -0.000608601100755223
$ws.Cells.Item(18, 5).Value = -2.02920667547382
$ws.Cells.Item(18, 6).Value = 0.101386015116349
$ws.Cells.Item(18, 7).Value = -0.165252558924571
$ws.Cells.Item(18, 8).Value = 0.04858091302945
$ws.Cells.Item(18, 9).Value = 0.0102791240611726
$ws.Cells.Item(18, 10).Value = 0.00236010511077498
$ws.Cells.Item(18, 11).Value = -0.0018486256997203
$ws.Cells.Item(19, 5).Value = -2.88701180933812
$ws.Cells.Item(19, 6).Value = 0.187647325679811
$ws.Cells.Item(19, 7).Value = 0.161475571525528
$ws.Cells.Item(19, 8).Value = 0.0969665892615014
$ws.Cells.Item(19, 9).Value = 0.035211518834785
$ws.Cells.Item(19, 10).Value = 0.00940251943300871
$ws.Cells.Item(19, 11).Value = -0.0137431723158013
$ws.Cells.Item(20, 5).Value = -2.80076969518453
$ws.Cells.Item(20, 6).Value = 0.199420834984694
$ws.Cells.Item(20, 7).Value = 0.300144465918335
$ws.Cells.Item(20, 8).Value = 0.135723151780683
$ws.Cells.Item(20, 9).Value = 0.0397686694259927
$ws.Cells.Item(20, 10).Value = 0.0184207739292822
$ws.Cells.Item(20, 11).Value = -0.0225109028945209
$ws.Cells.Item(21, 5).Value = -1.95308165608033
$ws.Cells.Item(21, 6).Value = 0.235764186065949
$ws.Cells.Item(21, 7).Value = -0.127807148941032
$ws.Cells.Item(21, 8).Value = 0.128157101154197
$ws.Cells.Item(21, 9).Value = 0.0555847514313396
$ws.Cells.Item(21, 10).Value = 0.016424242576247
$ws.Cells.Item(21, 11).Value = -0.014699123668472
$ws.Cells.Item(22, 5).Value = -2.70892943143692
$ws.Cells.Item(22, 6).Value = 0.395908632478748
$ws.Cells.Item(22, 7).Value = 0.242382907798554
$ws.Cells.Item(22, 8).Value = 0.219530970391322
$ws.Cells.Item(22, 9).Value = 0.156743645271192
$ws.Cells.Item(22, 10).Value = 0.0481938469609553
$ws.Cells.Item(22, 11).Value = -0.0662510125415888

$ws = $wb.Worksheets.Item("data2 lognormal")
$ws.Cells.Item(2, 5).Value = 2.03645638169679
$ws.Cells.Item(2, 6).Value = 0.152717764270313
$ws.Cells.Item(2, 7).Value = -1.12190166059483
$ws.Cells.Item(2, 8).Value = 0.0867906653867666
$ws.Cells.Item(2, 9).Value = 0.0233227155237229
$ws.Cells.Item(2, 10).Value = 0.00753261959827769
$ws.Cells.Item(2, 11).Value = -0.0107506639965146
$ws.Cells.Item(3, 5).Value = 2.0046720337807
$ws.Cells.Item(3, 6).Value = 0.150216926418198
$ws.Cells.Item(3, 7).Value = -0.992824174200785
$ws.Cells.Item(3, 8).Value = 0.0678852275553906
$ws.Cells.Item(3, 9).Value = 0.0225651249825304
$ws.Cells.Item(3, 10).Value = 0.00460840412024716
$ws.Cells.Item(3, 11).Value = -0.00802013418412126
$ws.Cells.Item(4, 5).Value = 2.51392827451678
$ws.Cells.Item(4, 6).Value = 0.201082049136534
$ws.Cells.Item(4, 7).Value = -1.04507703041139
$ws.Cells.Item(4, 8).Value = 0.11006605620096
$ws.Cells.Item(4, 9).Value = 0.0404339904849474
$ws.Cells.Item(4, 10).Value = 0.0121145367276329
$ws.Cells.Item(4, 11).Value = -0.0187037178595124
$ws.Cells.Item(5, 5).Value = 2.90225813456849
$ws.Cells.Item(5, 6).Value = 0.237300586670161
$ws.Cells.Item(5, 7).Value = -1.14817099226754
$ws.Cells.Item(5, 8).Value = 0.105467977729247
$ws.Cells.Item(5, 9).Value = 0.0563115684340025
$ws.Cells.Item(5, 10).Value = 0.0111234943262969
$ws.Cells.Item(5, 11).Value = -0.021248020735707
$ws.Cells.Item(6, 5).Value = 2.08586326373053
$ws.Cells.Item(6, 6).Value = 0.145302181619435
$ws.Cells.Item(6, 7).Value = -1.03120411829495
$ws.Cells.Item(6, 8).Value = 0.074283046751274
$ws.Cells.Item(6, 9).Value = 0.0211127239833672
$ws.Cells.Item(6, 10).Value = 0.00551797103465195
$ws.Cells.Item(6, 11).Value = -0.00869402226637252
$ws.Cells.Item(7, 5).Value = 2.37960527088054
$ws.Cells.Item(7, 6).Value = 0.179786460979122
$ws.Cells.Item(7, 7).Value = -1.07901445539707
$ws.Cells.Item(7, 8).Value = 0.10031066104669
$ws.Cells.Item(7, 9).Value = 0.0323231715513972
$ws.Cells.Item(7, 10).Value = 0.010062228719624
$ws.Cells.Item(7, 11).Value = -0.0162601111896603
$ws.Cells.Item(8, 5).Value = 1.52166620795309
$ws.Cells.Item(8, 6).Value = 0.333572687350795
$ws.Cells.Item(8, 7).Value = -0.773932916603465
$ws.Cells.Item(8, 8).Value = 0.135974053392742
$ws.Cells.Item(8, 9).Value = 0.111270737746431
$ws.Cells.Item(8, 10).Value = 0.0184889431960524
$ws.Cells.Item(8, 11).Value = -0.0303399210140507
$ws.Cells.Item(9, 5).Value = 2.48739055832804
$ws.Cells.Item(9, 6).Value = 0.423446230224185
$ws.Cells.Item(9, 7).Value = -1.07834058869686
$ws.Cells.Item(9, 8).Value = 0.159202180743191
$ws.Cells.Item(9, 9).Value = 0.179306709891073
$ws.Cells.Item(9, 10).Value = 0.0253453343533878
$ws.Cells.Item(9, 11).Value = -0.0627180539399422
$ws.Cells.Item(10, 5).Value = 1.49818030642701
$ws.Cells.Item(10, 6).Value = 0.331222567339514
$ws.Cells.Item(10, 7).Value = -1.01418267661033
$ws.Cells.Item(10, 8).Value = 0.147453826849644
$ws.Cells.Item(10, 9).Value = 0.109708389114979
$ws.Cells.Item(10, 10).Value = 0.0217426310526049
$ws.Cells.Item(10, 11).Value = -0.0418133861145634
$ws.Cells.Item(11, 5).Value = 2.01754446692552
$ws.Cells.Item(11, 6).Value = 0.317384281227654
$ws.Cells.Item(11, 7).Value = -0.984930211873869
$ws.Cells.Item(11, 8).Value = 0.101801443691546
$ws.Cells.Item(11, 9).Value = 0.100732781970394
$ws.Cells.Item(11, 10).Value = 0.010363533937683
$ws.Cells.Item(11, 11).Value = -0.0284881700528722
$ws.Cells.Item(12, 5).Value = 2.02696016167979
$ws.Cells.Item(12, 6).Value = 0.409698626967328
$ws.Cells.Item(12, 7).Value = -0.968032021305164
$ws.Cells.Item(12, 8).Value = 0.161971004464799
$ws.Cells.Item(12, 9).Value = 0.167852964938914
$ws.Cells.Item(12, 10).Value = 0.026234606287336
$ws.Cells.Item(12, 11).Value = -0.0587283707432745
$ws.Cells.Item(13, 5).Value = 2.88150849650292
$ws.Cells.Item(13, 6).Value = 0.698101806909137
$ws.Cells.Item(13, 7).Value = -1.04537780645883
$ws.Cells.Item(13, 8).Value = 0.221833270937624
$ws.Cells.Item(13, 9).Value = 0.487346132809803
$ws.Cells.Item(13, 10).Value = 0.0492100000948854
$ws.Cells.Item(13, 11).Value = -0.149652715148843
$ws.Cells.Item(14, 5).Value = 2.24237244057441
$ws.Cells.Item(14, 6).Value = 0.209726581516632
$ws.Cells.Item(14, 7).Value = -0.947676440017259
$ws.Cells.Item(14, 8).Value = 0.093789380646873
$ws.Cells.Item(14, 9).Value = 0.0439852389946526
$ws.Cells.Item(14, 10).Value = 0.00879644792212403
$ws.Cells.Item(14, 11).Value = -0.0161834686181826
$ws.Cells.Item(15, 5).Value = 2.46094374220508
$ws.Cells.Item(15, 6).Value = 0.264036571292589
$ws.Cells.Item(15, 7).Value = -1.01546611467525
$ws.Cells.Item(15, 8).Value = 0.109562714690895
$ws.Cells.Item(15, 9).Value = 0.0697153109799464
$ws.Cells.Item(15, 10).Value = 0.0120039884504384
$ws.Cells.Item(15, 11).Value = -0.0254842257230339
$ws.Cells.Item(16, 5).Value = 2.25444620048039
$ws.Cells.Item(16, 6).Value = 0.222341048830984
$ws.Cells.Item(16, 7).Value = -0.911033282332893
$ws.Cells.Item(16, 8).Value = 0.102491288079624
$ws.Cells.Item(16, 9).Value = 0.0494355419952621
$ws.Cells.Item(16, 10).Value = 0.0105044641322206
$ws.Cells.Item(16, 11).Value = -0.0192375999914305
$ws.Cells.Item(17, 5).Value = 0.996763855679523
$ws.Cells.Item(17, 6).Value = 0.0838183464811693
$ws.Cells.Item(17, 7).Value = -0.717453756337411
$ws.Cells.Item(17, 8).Value = 0.0396326427899414
$ws.Cells.Item(17, 9).Value = 0.00702551520683735
$ws.Cells.Item(17, 10).Value = 0.00157074637451509
$ws.Cells.Item(17, 11).Value = -0.000964577313751146
$ws.Cells.Item(18, 5).Value = 1.12639172344645
$ws.Cells.Item(18, 6).Value = 0.105907986747605
$ws.Cells.Item(18, 7).Value = -0.699866140753643
$ws.Cells.Item(18, 8).Value = 0.0462881852377329
$ws.Cells.Item(18, 9).Value = 0.0112165016569309
$ws.Cells.Item(18, 10).Value = 0.00214259609260268
$ws.Cells.Item(18, 11).Value = -0.00303414204430054
$ws.Cells.Item(19, 5).Value = 2.30823999869546
$ws.Cells.Item(19, 6).Value = 0.212100236552518
$ws.Cells.Item(19, 7).Value = -1.06317174036088
$ws.Cells.Item(19, 8).Value = 0.0883723508302685
$ws.Cells.Item(19, 9).Value = 0.0449865103456341
$ws.Cells.Item(19, 10).Value = 0.00780967239126806
$ws.Cells.Item(19, 11).Value = -0.0158401672908242
$ws.Cells.Item(20, 5).Value = 2.39272083727114
$ws.Cells.Item(20, 6).Value = 0.265967582835733
$ws.Cells.Item(20, 7).Value = -1.21578962659053
$ws.Cells.Item(20, 8).Value = 0.119222995723348
$ws.Cells.Item(20, 9).Value = 0.0707387551194827
$ws.Cells.Item(20, 10).Value = 0.0142141227092495
$ws.Cells.Item(20, 11).Value = -0.0289436015331549
$ws.Cells.Item(21, 5).Value = 1.05257362849532
$ws.Cells.Item(21, 6).Value = 0.288080901041177
$ws.Cells.Item(21, 7).Value = -0.717436925236619
$ws.Cells.Item(21, 8).Value = 0.124537607141585
$ws.Cells.Item(21, 9).Value = 0.0829906055446962
$ws.Cells.Item(21, 10).Value = 0.0155096155925518
$ws.Cells.Item(21, 11).Value = -0.0278016389585685
$ws.Cells.Item(22, 5).Value = 1.89257104125165
$ws.Cells.Item(22, 6).Value = 0.410670123606949
$ws.Cells.Item(22, 7).Value = -0.987198951461175
$ws.Cells.Item(22, 8).Value = 0.173663665314817
$ws.Cells.Item(22, 9).Value = 0.168649950423347
$ws.Cells.Item(22, 10).Value = 0.0301590686505767
$ws.Cells.Item(22, 11).Value = -0.0611261870872031

$ws = $wb.Worksheets.Item("data2 llogis")
$ws.Cells.Item(2, 5).Value = -1.777951348557
$ws.Cells.Item(2, 6).Value = 0.0825379900743899
$ws.Cells.Item(2, 7).Value = 0.679424176932772
$ws.Cells.Item(2, 8).Value = 0.0853173482110215
$ws.Cells.Item(2, 9).Value = 0.00681251980552008
$ws.Cells.Item(2, 10).Value = 0.00727904990576069
$ws.Cells.Item(2, 11).Value = 0.00220984086502586
$ws.Cells.Item(3, 5).Value = -1.99340653773269
$ws.Cells.Item(3, 6).Value = 0.0976327313995331
$ws.Cells.Item(3, 7).Value = 0.535253073568705
$ws.Cells.Item(3, 8).Value = 0.0774279575833975
$ws.Cells.Item(3, 9).Value = 0.00953215024053337
$ws.Cells.Item(3, 10).Value = 0.00599508861553641
$ws.Cells.Item(3, 11).Value = 0.0017118540801892
$ws.Cells.Item(4, 5).Value = -2.36474220234867
$ws.Cells.Item(4, 6).Value = 0.137854567634377
$ws.Cells.Item(4, 7).Value = 0.55764737964053
$ws.Cells.Item(4, 8).Value = 0.111678500072081
$ws.Cells.Item(4, 9).Value = 0.019003881817661
$ws.Cells.Item(4, 10).Value = 0.0124720873783497
$ws.Cells.Item(4, 11).Value = 0.00905882053789058
$ws.Cells.Item(5, 5).Value = -2.5105798322056
$ws.Cells.Item(5, 6).Value = 0.127186380028304
$ws.Cells.Item(5, 7).Value = 0.647579535658755
$ws.Cells.Item(5, 8).Value = 0.0977150205592001
$ws.Cells.Item(5, 9).Value = 0.0161763752647042
$ws.Cells.Item(5, 10).Value = 0.00954822524288489
$ws.Cells.Item(5, 11).Value = 0.00484636394239416
$ws.Cells.Item(6, 5).Value = -1.98994626729784
$ws.Cells.Item(6, 6).Value = 0.0920141316934637
$ws.Cells.Item(6, 7).Value = 0.582431980739211
$ws.Cells.Item(6, 8).Value = 0.0716241750542441
$ws.Cells.Item(6, 9).Value = 0.00846660043130209
$ws.Cells.Item(6, 10).Value = 0.00513002245220101
$ws.Cells.Item(6, 11).Value = 0.00154409368478559
$ws.Cells.Item(7, 5).Value = -2.12776853505119
$ws.Cells.Item(7, 6).Value = 0.0802505639146212
$ws.Cells.Item(7, 7).Value = 0.682229582156072
$ws.Cells.Item(7, 8).Value = 0.102850940330416
$ws.Cells.Item(7, 9).Value = 0.0064401530086147
$ws.Cells.Item(7, 10).Value = 0.0105783159268507
$ws.Cells.Item(7, 11).Value = 0.00421831753931453
$ws.Cells.Item(8, 5).Value = -1.9482535140164
$ws.Cells.Item(8, 6).Value = 0.320334372530899
$ws.Cells.Item(8, 7).Value = 0.261024068593365
$ws.Cells.Item(8, 8).Value = 0.189629506963794
$ws.Cells.Item(8, 9).Value = 0.102614110224765
$ws.Cells.Item(8, 10).Value = 0.0359593499113315
$ws.Cells.Item(8, 11).Value = 0.00923537545251059
$ws.Cells.Item(9, 5).Value = -2.3202820840382
$ws.Cells.Item(9, 6).Value = 0.12396398295321
$ws.Cells.Item(9, 7).Value = 0.686335169404983
$ws.Cells.Item(9, 8).Value = 0.16539975502107
$ws.Cells.Item(9, 9).Value = 0.0153670690696236
$ws.Cells.Item(9, 10).Value = 0.02735707896103
$ws.Cells.Item(9, 11).Value = -0.00232816030093871
$ws.Cells.Item(10, 5).Value = -1.51780137431937
$ws.Cells.Item(10, 6).Value = 0.174083053785454
$ws.Cells.Item(10, 7).Value = 0.580710717024211
$ws.Cells.Item(10, 8).Value = 0.149819611981135
$ws.Cells.Item(10, 9).Value = 0.0303049096152694
$ws.Cells.Item(10, 10).Value = 0.0224459161341779
$ws.Cells.Item(10, 11).Value = -0.0072368232351397
$ws.Cells.Item(11, 5).Value = -2.1068711084345
$ws.Cells.Item(11, 6).Value = 0.179161935426299
$ws.Cells.Item(11, 7).Value = 0.523639577340592
$ws.Cells.Item(11, 8).Value = 0.11890355424349
$ws.Cells.Item(11, 9).Value = 0.0320989991056973
$ws.Cells.Item(11, 10).Value = 0.0141380552117345
$ws.Cells.Item(11, 11).Value = -0.00954686235042087
$ws.Cells.Item(12, 5).Value = -2.12000443797963
$ws.Cells.Item(12, 6).Value = 0.195730672493874
$ws.Cells.Item(12, 7).Value = 0.506820763389294
$ws.Cells.Item(12, 8).Value = 0.170291094817232
$ws.Cells.Item(12, 9).Value = 0.0383104961549042
$ws.Cells.Item(12, 10).Value = 0.0289990569740515
$ws.Cells.Item(12, 11).Value = -0.00378819639480207
$ws.Cells.Item(13, 5).Value = -2.79988513554914
$ws.Cells.Item(13, 6).Value = 0.164643636088581
$ws.Cells.Item(13, 7).Value = 0.629808504017636
$ws.Cells.Item(13, 8).Value = 0.206972626700285
$ws.Cells.Item(13, 9).Value = 0.0271075269044691
$ws.Cells.Item(13, 10).Value = 0.0428376682032157
$ws.Cells.Item(13, 11).Value = -0.00816904122751264
$ws.Cells.Item(14, 5).Value = -2.33221811018217
$ws.Cells.Item(14, 6).Value = 0.13146167929006
$ws.Cells.Item(14, 7).Value = 0.487257647581758
$ws.Cells.Item(14, 8).Value = 0.110417509134753
$ws.Cells.Item(14, 9).Value = 0.0172821731217626
$ws.Cells.Item(14, 10).Value = 0.0121920263235233
$ws.Cells.Item(14, 11).Value = 0.00414207079880718
$ws.Cells.Item(15, 5).Value = -2.40113898532518
$ws.Cells.Item(15, 6).Value = 0.122189329723341
$ws.Cells.Item(15, 7).Value = 0.565964498875003
$ws.Cells.Item(15, 8).Value = 0.112695489584321
$ws.Cells.Item(15, 9).Value = 0.0149302322982393
$ws.Cells.Item(15, 10).Value = 0.0127002733726498
$ws.Cells.Item(15, 11).Value = 0.00300575235300295
$ws.Cells.Item(16, 5).Value = -2.4143783123601
$ws.Cells.Item(16, 6).Value = 0.137942025672175
$ws.Cells.Item(16, 7).Value = 0.462959154088631
$ws.Cells.Item(16, 8).Value = 0.128537936727382
$ws.Cells.Item(16, 9).Value = 0.0190280024465431
$ws.Cells.Item(16, 10).Value = 0.0165220011781323
$ws.Cells.Item(16, 11).Value = 0.00808717390768899
$ws.Cells.Item(17, 5).Value = -1.35404457784318
$ws.Cells.Item(17, 6).Value = 0.128536016420039
$ws.Cells.Item(17, 7).Value = 0.175614504858684
$ws.Cells.Item(17, 8).Value = 0.0605189546240113
$ws.Cells.Item(17, 9).Value = 0.0165215075171325
$ws.Cells.Item(17, 10).Value = 0.00366254386878314
$ws.Cells.Item(17, 11).Value = 0.0023704957571454
$ws.Cells.Item(18, 5).Value = -1.61323348320391
$ws.Cells.Item(18, 6).Value = 0.121792974479335
$ws.Cells.Item(18, 7).Value = 0.180586724304846
$ws.Cells.Item(18, 8).Value = 0.0673717463977423
$ws.Cells.Item(18, 9).Value = 0.0148335286325239
$ws.Cells.Item(18, 10).Value = 0.0045389522126817
$ws.Cells.Item(18, 11).Value = 0.00101856287708708
$ws.Cells.Item(19, 5).Value = -2.16821355133828
$ws.Cells.Item(19, 6).Value = 0.108395064222108
$ws.Cells.Item(19, 7).Value = 0.599206094165685
$ws.Cells.Item(19, 8).Value = 0.0902889678466502
$ws.Cells.Item(19, 9).Value = 0.0117494899477149
$ws.Cells.Item(19, 10).Value = 0.00815209771481344
$ws.Cells.Item(19, 11).Value = 0.000589220943448199
$ws.Cells.Item(20, 5).Value = -1.99064811632014
$ws.Cells.Item(20, 6).Value = 0.0821751582283421
$ws.Cells.Item(20, 7).Value = 0.785020862239596
$ws.Cells.Item(20, 8).Value = 0.106565176953077
$ws.Cells.Item(20, 9).Value = 0.00675275662985306
$ws.Cells.Item(20, 10).Value = 0.0113561369390407
$ws.Cells.Item(20, 11).Value = -0.00158573335259087
$ws.Cells.Item(21, 5).Value = -1.50596150897498
$ws.Cells.Item(21, 6).Value = 0.249752256684582
$ws.Cells.Item(21, 7).Value = 0.222267720042517
$ws.Cells.Item(21, 8).Value = 0.166463011139801
$ws.Cells.Item(21, 9).Value = 0.0623761897190414
$ws.Cells.Item(21, 10).Value = 0.0277099340777296
$ws.Cells.Item(21, 11).Value = -0.00515070683471067
$ws.Cells.Item(22, 5).Value = -1.94606176905585
$ws.Cells.Item(22, 6).Value = 0.208222229400172
$ws.Cells.Item(22, 7).Value = 0.536281630291043
$ws.Cells.Item(22, 8).Value = 0.18904560480009
$ws.Cells.Item(22, 9).Value = 0.0433564968163777
$ws.Cells.Item(22, 10).Value = 0.0357382406942316
$ws.Cells.Item(22, 11).Value = -0.00355768253319537

$ws = $wb.Worksheets.Item("data2 gompertz")
$ws.Cells.Item(2, 5).Value = -2.08140174659871
$ws.Cells.Item(2, 6).Value = 0.101464536652
$ws.Cells.Item(2, 7).Value = -0.0180406304756775
$ws.Cells.Item(2, 8).Value = 0.012319060510442
$ws.Cells.Item(2, 9).Value = 0.0102950521980051
$ws.Cells.Item(2, 10).Value = 0.000151759251859932
$ws.Cells.Item(2, 11).Value = -0.000601428257233394
$ws.Cells.Item(3, 5).Value = -2.3636120790641
$ws.Cells.Item(3, 6).Value = 0.126087808930259
$ws.Cells.Item(3, 7).Value = -0.0125577421591243
$ws.Cells.Item(3, 8).Value = 0.00856252069692879
$ws.Cells.Item(3, 9).Value = 0.0158981355608334
$ws.Cells.Item(3, 10).Value = 0.0000733167606853339
$ws.Cells.Item(3, 11).Value = -0.000677055669639301
$ws.Cells.Item(4, 5).Value = -2.70193626478097
$ws.Cells.Item(4, 6).Value = 0.156444834207228
$ws.Cells.Item(4, 7).Value = -0.0190879804587354
$ws.Cells.Item(4, 8).Value = 0.019550634465166
$ws.Cells.Item(4, 9).Value = 0.0244749861501269
$ws.Cells.Item(4, 10).Value = 0.000382227307990537
$ws.Cells.Item(4, 11).Value = -0.00172192507099017
$ws.Cells.Item(5, 5).Value = -3.06419331756653
$ws.Cells.Item(5, 6).Value = 0.185570645347155
$ws.Cells.Item(5, 7).Value = 0.0121285462350868
$ws.Cells.Item(5, 8).Value = 0.0163362747356503
$ws.Cells.Item(5, 9).Value = 0.0344364644145594
$ws.Cells.Item(5, 10).Value = 0.000266873872238646
$ws.Cells.Item(5, 11).Value = -0.00195349293902929
$ws.Cells.Item(6, 5).Value = -2.23135235698007
$ws.Cells.Item(6, 6).Value = 0.101177196866545
$ws.Cells.Item(6, 7).Value = -0.0235283723645498
$ws.Cells.Item(6, 8).Value = 0.00996964691214282
$ws.Cells.Item(6, 9).Value = 0.0102368251657716
$ws.Cells.Item(6, 10).Value = 0.0000993938595527989
$ws.Cells.Item(6, 11).Value = -0.000595977887148311
$ws.Cells.Item(7, 5).Value = -2.30353901599517
$ws.Cells.Item(7, 6).Value = 0.0849990247310865
$ws.Cells.Item(7, 7).Value = -0.0295289955002821
$ws.Cells.Item(7, 8).Value = 0.00740228834920616
$ws.Cells.Item(7, 9).Value = 0.00722483420523586
$ws.Cells.Item(7, 10).Value = 0.0000547938728047933
$ws.Cells.Item(7, 11).Value = -0.000207469315286281
$ws.Cells.Item(8, 5).Value = -2.2271083115886
$ws.Cells.Item(8, 6).Value = 0.34245155985288
$ws.Cells.Item(8, 7).Value = -0.0377432470081011
$ws.Cells.Item(8, 8).Value = 0.030834685943283
$ws.Cells.Item(8, 9).Value = 0.117273070845671
$ws.Cells.Item(8, 10).Value = 0.000950777857220891
$ws.Cells.Item(8, 11).Value = -0.00652736343257616
$ws.Cells.Item(9, 5).Value = -2.65972894717063
$ws.Cells.Item(9, 6).Value = 0.189263107546043
$ws.Cells.Item(9, 7).Value = -0.00246882176944569
$ws.Cells.Item(9, 8).Value = 0.0180428086912718
$ws.Cells.Item(9, 9).Value = 0.0358205238779849
$ws.Cells.Item(9, 10).Value = 0.000325542945469833
$ws.Cells.Item(9, 11).Value = -0.00230345389205733
$ws.Cells.Item(10, 5).Value = -2.02381715538323
$ws.Cells.Item(10, 6).Value = 0.255778640267797
$ws.Cells.Item(10, 7).Value = 0.0255597996111659
$ws.Cells.Item(10, 8).Value = 0.0404981984260776
$ws.Cells.Item(10, 9).Value = 0.0654227128172429
$ws.Cells.Item(10, 10).Value = 0.00164010407575795
$ws.Cells.Item(10, 11).Value = -0.00711036413780034
$ws.Cells.Item(11, 5).Value = -2.67178985962201
$ws.Cells.Item(11, 6).Value = 0.230673825428547
$ws.Cells.Item(11, 7).Value = 0.0177692748011179
$ws.Cells.Item(11, 8).Value = 0.0167409973495209
$ws.Cells.Item(11, 9).Value = 0.0532104137378396
$ws.Cells.Item(11, 10).Value = 0.000280260992256664
$ws.Cells.Item(11, 11).Value = -0.00252369250313782
$ws.Cells.Item(12, 5).Value = -2.61749017114501
$ws.Cells.Item(12, 6).Value = 0.28894923044443
$ws.Cells.Item(12, 7).Value = 0.0123792653404658
$ws.Cells.Item(12, 8).Value = 0.0306305115929061
$ws.Cells.Item(12, 9).Value = 0.0834916577744283
$ws.Cells.Item(12, 10).Value = 0.000938228240443155
$ws.Cells.Item(12, 11).Value = -0.0063487416996439
$ws.Cells.Item(13, 5).Value = -3.46876419652843
$ws.Cells.Item(13, 6).Value = 0.299780689171852
$ws.Cells.Item(13, 7).Value = 0.0242905260330743
$ws.Cells.Item(13, 8).Value = 0.0184675218673625
$ws.Cells.Item(13, 9).Value = 0.0898684616003503
$ws.Cells.Item(13, 10).Value = 0.000341049363921511
$ws.Cells.Item(13, 11).Value = -0.0041447436937204
$ws.Cells.Item(14, 5).Value = -2.53774537520061
$ws.Cells.Item(14, 6).Value = 0.141652278771847
$ws.Cells.Item(14, 7).Value = -0.0256502815535124
$ws.Cells.Item(14, 8).Value = 0.0115962452167281
$ws.Cells.Item(14, 9).Value = 0.0200653680812571
$ws.Cells.Item(14, 10).Value = 0.00013447290312649
$ws.Cells.Item(14, 11).Value = -0.000884315954813601
$ws.Cells.Item(15, 5).Value = -2.64863060347219
$ws.Cells.Item(15, 6).Value = 0.140374272993617
$ws.Cells.Item(15, 7).Value = -0.0168393407736236
$ws.Cells.Item(15, 8).Value = 0.0119133252065979
$ws.Cells.Item(15, 9).Value = 0.0197049365184865
$ws.Cells.Item(15, 10).Value = 0.000141927317478162
$ws.Cells.Item(15, 11).Value = -0.000937601013814208
$ws.Cells.Item(16, 5).Value = -2.51123869942539
$ws.Cells.Item(16, 6).Value = 0.121563063431743
$ws.Cells.Item(16, 7).Value = -0.0369299932610897
$ws.Cells.Item(16, 8).Value = 0.010499509130997
$ws.Cells.Item(16, 9).Value = 0.01477757839091
$ws.Cells.Item(16, 10).Value = 0.000110239691991889
$ws.Cells.Item(16, 11).Value = -0.000462489919387872
$ws.Cells.Item(17, 5).Value = -1.56961926885366
$ws.Cells.Item(17, 6).Value = 0.130599670361597
$ws.Cells.Item(17, 7).Value = -0.0906814699721602
$ws.Cells.Item(17, 8).Value = 0.0202486490166893
$ws.Cells.Item(17, 9).Value = 0.0170562738985578
$ws.Cells.Item(17, 10).Value = 0.000410007787001074
$ws.Cells.Item(17, 11).Value = -0.00191837604648169
$ws.Cells.Item(18, 5).Value = -1.84249044736258
$ws.Cells.Item(18, 6).Value = 0.127874783169237
$ws.Cells.Item(18, 7).Value = -0.0639210532659732
$ws.Cells.Item(18, 8).Value = 0.0168349599084398
$ws.Cells.Item(18, 9).Value = 0.0163519601705793
$ws.Cells.Item(18, 10).Value = 0.000283415875118776
$ws.Cells.Item(18, 11).Value = -0.00153860905264522
$ws.Cells.Item(19, 5).Value = -2.56149356438312
$ws.Cells.Item(19, 6).Value = 0.14962558407627
$ws.Cells.Item(19, 7).Value = -0.00129002285038163
$ws.Cells.Item(19, 8).Value = 0.0138376768307981
$ws.Cells.Item(19, 9).Value = 0.0223878154101648
$ws.Cells.Item(19, 10).Value = 0.000191481300073606
$ws.Cells.Item(19, 11).Value = -0.00137946959787476
$ws.Cells.Item(20, 5).Value = -2.34740021745688
$ws.Cells.Item(20, 6).Value = 0.126118366223994
$ws.Cells.Item(20, 7).Value = 0.00794657992825708
$ws.Cells.Item(20, 8).Value = 0.0201871014722952
$ws.Cells.Item(20, 9).Value = 0.0159058422990094
$ws.Cells.Item(20, 10).Value = 0.000407519065852742
$ws.Cells.Item(20, 11).Value = -0.0017874632243433
$ws.Cells.Item(21, 5).Value = -1.8334625702322
$ws.Cells.Item(21, 6).Value = 0.270318049401206
$ws.Cells.Item(21, 7).Value = -0.0444484446995318
$ws.Cells.Item(21, 8).Value = 0.0368542173565487
$ws.Cells.Item(21, 9).Value = 0.0730718478320728
$ws.Cells.Item(21, 10).Value = 0.00135823333696374
$ws.Cells.Item(21, 11).Value = -0.00680553207578169
$ws.Cells.Item(22, 5).Value = -2.50741289770736
$ws.Cells.Item(22, 6).Value = 0.318208850704761
$ws.Cells.Item(22, 7).Value = 0.0248428547262417
$ws.Cells.Item(22, 8).Value = 0.047688533074933
$ws.Cells.Item(22, 9).Value = 0.101256872666845
$ws.Cells.Item(22, 10).Value = 0.00227419618683898
$ws.Cells.Item(22, 11).Value = -0.0101235506567224
